# Apply the changes described by the diff:
#  - Rename Sheet1 -> AddPriceAgreement
#  - Populate AddPriceAgreement (sheet2) with Role/Location/selectUOM headers
#    and REQUESTOR/XEEVA -MJ/CU-CUBIC data row
#  - Make AddPriceAgreement the active sheet/tab, with E4 selected
#  - On the Data sheet (sheet1), change the selection to A1:B2 (no longer the
#    active tab)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Data"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1" -> "AddPriceAgreement"

# Rename the second worksheet
$ws2.Name = "AddPriceAgreement"

# Fill in the new table. Set the column-3 values first (CU-CUBIC before
# selectUOM) so the shared-string table receives entries in the same order
# as the target workbook.
$ws2.Cells.Item(1, 1).Value = "Role"
$ws2.Cells.Item(2, 1).Value = "REQUESTOR"
$ws2.Cells.Item(1, 2).Value = "Location"
$ws2.Cells.Item(2, 2).Value = "XEEVA -MJ"
$ws2.Cells.Item(2, 3).Value = "CU-CUBIC"
$ws2.Cells.Item(1, 3).Value = "selectUOM"

# Size the new columns to fit their content
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()

# Select A1:B2 on the Data sheet
$ws1.Range("A1:B2").Select() | Out-Null

# Make AddPriceAgreement the active sheet/tab and select E4 on it
$ws2.Activate() | Out-Null
$ws2.Range("E4").Select() | Out-Null
